$wb = $excel.ActiveWorkbook
$wsProjects = $wb.Worksheets.Item("#Projects")
$wsPeople   = $wb.Worksheets.Item("#People")

# --- "#Projects" sheet: add a new "member" column (H) that shows, per
#     project, the comma-separated list of person-ids working on it.
#     This replaces the old "[Members]" section label in A10 (that table
#     still exists below, but its header label is cleared since the
#     relationship is now also visible as a single multi-value column).

$wsProjects.Range("A10").Value = ""

$wsProjects.Range("H1").Value = "member"
$wsProjects.Range("H2").Value = "[Person,]"
$wsProjects.Range("H3").Value = "p10001,p10002,p10003"
$wsProjects.Range("H4").Value = "p10001,p10003,p10004,p10005"
$wsProjects.Range("H5").Value = "p10002,p10003,p10005"
$wsProjects.Range("H8").Value = "p10006,p10009,p10010"

# widen the new column to fit its content
$wsProjects.Columns.Item(8).ColumnWidth = 32.75

# --- "#People" sheet: widen the two trailing helper columns (F, G)
$wsPeople.Columns.Item(6).ColumnWidth = 14.25
$wsPeople.Columns.Item(7).ColumnWidth = 18.25

# --- selections / active sheet: "#Projects" becomes the active tab,
#     with A21 selected; "#People" keeps G11 selected for when the user
#     switches back to it.
$wsPeople.Range("G11").Select()
$wsProjects.Activate()
$wsProjects.Range("A21").Select()
